# Edit: "Show clear unauthorized message for RBAC and merchants pages"
#
# The status table in this doc is built from paragraphs whose single <w:r>
# run holds several <w:t> pieces separated by literal <w:tab/> elements.
# Word's Find/Replace (and Range.Text / InsertBefore / InsertAfter) in this
# host re-serializes the whole owning run as one <w:t>, collapsing the
# <w:tab/> run breaks into literal tab characters -- so for the two rows
# that change we rebuild the paragraph body with Range.InsertXML, which
# accepts a literal WordprocessingML fragment and preserves <w:tab/> as a
# distinct sibling element, matching how unedited rows remain in the file.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Row 1: "Platform RBAC Admin UI" -> "RBAC/Merchants Unauthorized UX" ---
$row1 = '<w:p ' + $wNs + '>' +
          '<w:pPr/>' +
          '<w:r>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/>' +
              '<w:sz w:val="24"/>' +
              '<w:sz-cs w:val="24"/>' +
            '</w:rPr>' +
            '<w:t xml:space="preserve">RBAC/Merchants Unauthorized UX</w:t>' +
            '<w:tab/>' +
            '<w:t xml:space="preserve">Developed</w:t>' +
            '<w:tab/>' +
            '<w:t xml:space="preserve">Added explicit unauthorized handling for 401/403 in Platform RBAC and Merchants admin pages; UI now shows: &quot;You are not authorized.&quot; instead of generic API errors.</w:t>' +
            '<w:tab/>' +
            '<w:t xml:space="preserve">None.</w:t>' +
          '</w:r>' +
        '</w:p>'

# --- Row 2: "Merchants Management UI" -> "Platform RBAC + Merchants UI" ---
$row2 = '<w:p ' + $wNs + '>' +
          '<w:pPr/>' +
          '<w:r>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/>' +
              '<w:sz w:val="24"/>' +
              '<w:sz-cs w:val="24"/>' +
            '</w:rPr>' +
            '<w:t xml:space="preserve">Platform RBAC + Merchants UI</w:t>' +
            '<w:tab/>' +
            '<w:t xml:space="preserve">Developed</w:t>' +
            '<w:tab/>' +
            '<w:t xml:space="preserve">Pages remain wired to real APIs and admin routes.</w:t>' +
            '<w:tab/>' +
            '<w:t xml:space="preserve">Role seeding still required for access.</w:t>' +
          '</w:r>' +
        '</w:p>'

$d.Paragraphs.Item(5).Range.InsertXML($row1)
$d.Paragraphs.Item(6).Range.InsertXML($row2)

# --- Drop the two rows that no longer apply: "Admin Navigation Wiring" and
#     "Single Register Entry" (now still paragraphs 7 and 8). Delete the
#     higher-indexed one first so the other index stays valid. ---
$d.Paragraphs.Item(8).Range.Delete()
$d.Paragraphs.Item(7).Range.Delete()

# --- Git state footer lines. They used to be paragraphs 11 and 12; the two
#     Range.Delete() calls above each remove a whole paragraph (mark
#     included), so everything after shifts up by 2 -> now paragraphs 9 and
#     10. Rebuilt the same way so xml:space="preserve" on the lone <w:t>
#     survives untouched, same as the sibling paragraphs around them. ---
$gitLine1 = '<w:p ' + $wNs + '>' +
              '<w:pPr/>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/>' +
                  '<w:sz w:val="24"/>' +
                  '<w:sz-cs w:val="24"/>' +
                '</w:rPr>' +
                '<w:t xml:space="preserve">- Last pushed commit: 61bee7a</w:t>' +
              '</w:r>' +
            '</w:p>'

$gitLine2 = '<w:p ' + $wNs + '>' +
              '<w:pPr/>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/>' +
                  '<w:sz w:val="24"/>' +
                  '<w:sz-cs w:val="24"/>' +
                '</w:rPr>' +
                '<w:t xml:space="preserve">- Current unauthorized-message UX update is local and pending push.</w:t>' +
              '</w:r>' +
            '</w:p>'

$d.Paragraphs.Item(9).Range.InsertXML($gitLine1)
$d.Paragraphs.Item(10).Range.InsertXML($gitLine2)
